$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 12:12"

# Row 4
$ws.Range("E4").Value = 2169006
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 153848

# Row 6
$ws.Range("B6").Value = 1588129
$ws.Range("C6").Value = 3745
$ws.Range("D6").Value = 1022606
$ws.Range("E6").Value = 530487
$ws.Range("G6").Value = 33
$ws.Range("H6").Value = 35036

# Row 13
$ws.Range("A13").Value = "Iran"
$ws.Range("B13").Value = 301530
$ws.Range("C13").Value = 2621
$ws.Range("D13").Value = 261200
$ws.Range("E13").Value = 23761
$ws.Range("G13").Value = 226
$ws.Range("H13").Value = 16569

# Row 14
$ws.Range("A14").Value = "Reino Unido"
$ws.Range("B14").Value = 301455
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("H14").Value = 45961

# Row 34
$ws.Range("B34").Value = 79159
$ws.Range("C34").Value = 590
$ws.Range("D34").Value = 61421
$ws.Range("E34").Value = 17317
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 421

# Row 37
$ws.Range("B37").Value = 68769
$ws.Range("C37").Value = 470
$ws.Range("D37").Value = 35516
$ws.Range("E37").Value = 32756
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 497

# Row 43
$ws.Range("B43").Value = 60223
$ws.Range("C43").Value = 302
$ws.Range("D43").Value = 53626
$ws.Range("E43").Value = 6248
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 349

# Row 47
$ws.Range("B47").Value = 49591
$ws.Range("C47").Value = 1356
$ws.Range("D47").Value = 26609
$ws.Range("E47").Value = 20678
$ws.Range("G47").Value = 35
$ws.Range("H47").Value = 2304

# Row 52
$ws.Range("E52").Value = 3246
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 145

# Row 54
$ws.Range("B54").Value = 36542
$ws.Range("C54").Value = 71
$ws.Range("D54").Value = 25471
$ws.Range("E54").Value = 9800

# Row 57
$ws.Range("B57").Value = 35022
$ws.Range("C57").Value = 220
$ws.Range("E57").Value = 1942
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 1980

# Row 74
$ws.Range("E74").Value = 4286
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 378

# Row 87
$ws.Range("A87").Value = "Consejo Danes para los Refugiados"
$ws.Range("B87").Value = 9010
$ws.Range("C87").Value = 79
$ws.Range("D87").Value = 6437
$ws.Range("E87").Value = 2358
$ws.Range("G87").Value = 5
$ws.Range("H87").Value = 215

# Row 88
$ws.Range("A88").Value = "Malasia"
$ws.Range("B88").Value = 8964
$ws.Range("C88").Value = 8
$ws.Range("D88").Value = 8617
$ws.Range("E88").Value = 223
$ws.Range("H88").Value = 124

# Row 90
$ws.Range("B90").Value = 7423
$ws.Range("C90").Value = 9
$ws.Range("E90").Value = 144

# Row 113
$ws.Range("A113").Value = "Hong Kong"
$ws.Range("B113").Value = 3152
$ws.Range("C113").Value = 149
$ws.Range("D113").Value = 1660
$ws.Range("E113").Value = 1468
$ws.Range("H113").Value = 24

# Row 114
$ws.Range("A114").Value = "Guinea Ecuatorial"
$ws.Range("B114").Value = 3071
$ws.Range("D114").Value = 842
$ws.Range("E114").Value = 2178
$ws.Range("H114").Value = 51

# Row 115
$ws.Range("A115").Value = "Montenegro"
$ws.Range("B115").Value = 3016
$ws.Range("D115").Value = 931
$ws.Range("E115").Value = 2038
$ws.Range("H115").Value = 47

# Row 118
$ws.Range("B118").Value = 2811
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 2333
$ws.Range("E118").Value = 467

# Row 125
$ws.Range("B125").Value = 2139
$ws.Range("C125").Value = 24
$ws.Range("D125").Value = 1780
$ws.Range("E125").Value = 242

Write-Host "edits applied"
